$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add new columns E:H with values 4,5,6,7 ---
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7

# Copy the formatting (bold font, border, centered alignment) from the
# existing header cell A1 onto the new header cells E1:H1.
$ws.Range("A1").Copy()
$ws.Range("E1:H1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# --- Updated statistics in columns C (std dev) and D (variance) ---
$ws.Range("C2").Value = 20.49631717546763
$ws.Range("D2").Value = 420.0990177573692

$ws.Range("C3").Value = 4.593024078404592
$ws.Range("D3").Value = 21.09587018480435

$ws.Range("C4").Value = 2.070909161059325
$ws.Range("D4").Value = 4.288664753359439

$ws.Range("C5").Value = 7.780698595839763
$ws.Range("D5").Value = 60.53927063930286

$ws.Range("C6").Value = 0.4933567175747401
$ws.Range("D6").Value = 0.2434008507761219

$ws.Range("C7").Value = 9.817534115584078
$ws.Range("D7").Value = 96.38397611065723

$ws.Range("C8").Value = 4.213680226897766
$ws.Range("D8").Value = 17.75510105454921

$ws.Range("C9").Value = 24.48105869165857
$ws.Range("D9").Value = 599.3222346644318

$ws.Range("C10").Value = 14.60895644455249
$ws.Range("D10").Value = 213.4216083988318

# --- New columns E (min), F (max), G (?) and H (range) for each data row ---
$ws.Range("E2").Value = 134
$ws.Range("F2").Value = 218
$ws.Range("G2").Value = 101
$ws.Range("H2").Value = 117

$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 31.2
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 31.2

$ws.Range("E4").Value = 4.34
$ws.Range("F4").Value = 15.33
$ws.Range("G4").Value = 0.98
$ws.Range("H4").Value = 14.35

$ws.Range("E5").Value = 26.115
$ws.Range("F5").Value = 42.49
$ws.Range("G5").Value = 6.74
$ws.Range("H5").Value = 35.75

$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1

$ws.Range("E7").Value = 53
$ws.Range("F7").Value = 78
$ws.Range("G7").Value = 13
$ws.Range("H7").Value = 65

$ws.Range("E8").Value = 25.805
$ws.Range("F8").Value = 46.58
$ws.Range("G8").Value = 14.7
$ws.Range("H8").Value = 31.88

$ws.Range("E9").Value = 7.51
$ws.Range("F9").Value = 147.19
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 147.19

$ws.Range("E10").Value = 45
$ws.Range("F10").Value = 64
$ws.Range("G10").Value = 15
$ws.Range("H10").Value = 49
